$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "Overview" — row 3 corresponds to b.md. Handoff report generated,
# so both language columns now show "Ready for handoff" and the overall
# "Latest HO Xliff Generate Date" is refreshed.
# ----------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-16 00:34:04"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c7fce62af9f447cd954a1cc2aa568dfbfbfc812e/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/282ea9fc621281f98205116040ba6ed46f51c8d0/e2e/b.md."

# ----------------------------------------------------------------------
# Sheet "zh-cn"
# ----------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# a.md (row 2) is now confirmed ready for handoff too.
$wsZhCn.Range("C2").Value = "Ready for handoff"

# b.md (row 3) got a fresh handoff package generated for it.
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-16 00:33:57"
$wsZhCn.Range("P3").Value = $errorDetail

# Error Detail column needs to be widened to fit the new message.
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

# ----------------------------------------------------------------------
# Sheet "de-de"
# ----------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

# a.md (row 2) is now confirmed ready for handoff too.
$wsDeDe.Range("C2").Value = "Ready for handoff"

# b.md (row 3) got a fresh handoff package generated for it.
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-16 00:34:04"
$wsDeDe.Range("P3").Value = $errorDetail

# Error Detail column needs to be widened to fit the new message.
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
